# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 76
$ws1.Range("F6").Value = 5314
$ws1.Range("F8").Value = 876
$ws1.Range("F10").Value = 2360
$ws1.Range("F13").Value = 2208

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 76
$ws4.Range("F6").Value = 5314
$ws4.Range("F10").Value = 876
$ws4.Range("F12").Value = 2360
$ws4.Range("F16").Value = 2208
